$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Overview" sheet: row 3 corresponds to b.md. Mark it ready for handoff and
# refresh the "Latest HO Xliff Generate Date" timestamp.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-25 12:39:20"

# ---------------------------------------------------------------------------
# "zh-cn" sheet: row 3 (b.md) has a new handoff file generated, which flips
# its status, clears the "duplicate content" flag, and records an error
# about the handback file being stale relative to the new handoff.
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-25 12:39:15"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6a7aad33685df75d37f7531af785c3ef277db43c/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b0f481c3d46193f7f41249cd65d7f48d496c1048/e2e/b.md."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.15

# ---------------------------------------------------------------------------
# "de-de" sheet: same update as zh-cn, but for the de-de handoff file.
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-25 12:39:20"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6a7aad33685df75d37f7531af785c3ef277db43c/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b0f481c3d46193f7f41249cd65d7f48d496c1048/e2e/b.md."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.15
